# Apply the "mapping" sheet changes:
#  - rename method "seqImmuCC" -> "seqimmucc" (rows 241-250, col A)
#  - rename method "DCQ" -> "dcq" (rows 251-269, col A), and pick up the
#    same formatting on col B that the rest of that block (e.g. C269) uses
#  - append a new "base" method block (rows 270-288) that duplicates the
#    method_cell_type / cell_type columns of the "dcq" block

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mapping")
$ws.Activate()

# --- rename seqImmuCC -> seqimmucc -----------------------------------
$ws.Range("A241:A250").Value = "seqimmucc"

# --- rename DCQ -> dcq -------------------------------------------------
$ws.Range("A251:A269").Value = "dcq"

# Column B in the old "DCQ" block had no explicit formatting; the edited
# workbook gives it the same style already used by C269 (Arial Unicode MS
# 10pt). Copy that formatting across so the new cells pick up the same
# style index instead of minting a new one.
$ws.Range("C269").Copy()
$ws.Range("B251:B269").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- add the new "base" block (duplicate of the dcq rows 251-269) -----
$baseData = @(
    @("B_cells", "B cell"),
    @("Dev_B_cells", "B cell naive"),
    @("Stromal", "stroma score"),
    @("Dendritic_cells", "Myeloid dendritic cell"),
    @("plasmacytoid_DCs", "Plasmacytoid dendritic cell"),
    @("Granulocytes", "Granulocyte-monocyte progenitor"),
    @("Macrophages", "Macrophage"),
    @("Stem", "Stem cell"),
    @("Monocytes", "Monocyte"),
    @("NK_cells", "NK cell"),
    @("NKT_cells", "T cell NK"),
    @("Thymocytes", "T cell immature"),
    @("CD4_Naive", "T cell CD4+ naive"),
    @("T-reg", "T cell regulatory (Tregs)"),
    @("CD4_Memory", "T cell CD4+ memory "),
    @("CD8_Eff", "T cell CD8+ effector memory"),
    @("CD8_Memory", "T cell CD8+ memory"),
    @("CD8_Naive", "T cell CD8+ naive"),
    @("T_cells_gamma_delta ", "T cell gamma delta")
)

$row = 270
foreach ($pair in $baseData) {
    $ws.Cells.Item($row, 1).Value = "base"
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
    $row = $row + 1
}

# --- restore the on-screen selection to roughly where the author left it
$ws.Range("A290").Select()

Write-Output "done"
